$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new data row (row 46) replicating the pattern of existing rows,
# for "Feria Lagunitas de Puerto Montt" / Granada.
$row = 46

$ws.Cells.Item($row, 1).Value = 4
$ws.Cells.Item($row, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($row, 3).Value = "Los Lagos"

$ws.Cells.Item($row, 4).Value = 45041
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item(45, 4).NumberFormat

$ws.Cells.Item($row, 5).Value = 10
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100104
$ws.Cells.Item($row, 8).Value = "Frutos de pepita"
$ws.Cells.Item($row, 9).Value = 100104001
$ws.Cells.Item($row, 10).Value = "Granada"
$ws.Cells.Item($row, 11).Value = "Wonderfull"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 200
$ws.Cells.Item($row, 14).Value = 16000
$ws.Cells.Item($row, 15).Value = 17000
$ws.Cells.Item($row, 16).Value = 16500
$ws.Cells.Item($row, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item($row, 18).Value = "Provincia de Limarí"
$ws.Cells.Item($row, 19).Value = 917
$ws.Cells.Item($row, 20).Value = 18
